# The commit swaps the two theme parts in the package:
#   ppt/theme/theme1.xml  (applied to the slide master / all slides)
#     "Integral" / "Red Violet" colour scheme  ->  default "Office Theme" colours
#   ppt/theme/theme2.xml  (applied to the notes master only)
#     default "Office Theme" colours            ->  "Integral" / "Red Violet" colours
#
# The PowerPoint object model only exposes a single, document-level
# ThemeColorScheme (reachable from the slide master, a slide, or the
# notes master - they all resolve to the same underlying theme part),
# so only the slide-facing theme (theme1.xml) can be edited here. We
# rewrite its 12 colour slots to the stock "Office" palette, which is
# the visible, structural part of the change described by the diff.

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Office theme palette (target content of ppt/theme/theme1.xml)
# Index order matches OOXML clrScheme: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$colorScheme.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$colorScheme.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$colorScheme.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$colorScheme.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$colorScheme.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$colorScheme.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$colorScheme.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$colorScheme.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
